# Move example_sds_dataset to resources: update sample rows and drop the
# now-redundant rows 6-8 (their sam-1/sam-2/sam-3 data has been folded into
# rows 2-5 alongside the updated sub-1/sub-2 pairing).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last three data rows (6-8); rows 2-5 below are rewritten with
# the merged values, so the tail rows are no longer needed.
$ws.Rows("6:8").Delete()

# Row 2: sam-1 / sub-1
$ws.Range("A2").Value = "sam-1"
$ws.Range("B2").Value = "sub-1"
$ws.Range("C2").Value = "Not Defined"
$ws.Range("E2").Value = "Experimental"
$ws.Range("F2").Value = "DCE-MRI Contrast Image sam-1"
$ws.Range("G2").Value = "Breast"
$ws.Range("H2").Value = "Not Defined"
$ws.Range("I2").Value = "Not Defined"

# Row 3: sam-2 / sub-1
$ws.Range("A3").Value = "sam-2"
$ws.Range("B3").Value = "sub-1"
$ws.Range("C3").Value = "Not Defined"
$ws.Range("E3").Value = "Experimental"
$ws.Range("F3").Value = "DCE-MRI Contrast Image sam-2"
$ws.Range("G3").Value = "Breast"
$ws.Range("H3").Value = "Not Defined"
$ws.Range("I3").Value = "Not Defined"

# Row 4: sam-1 / sub-2
$ws.Range("A4").Value = "sam-1"
$ws.Range("B4").Value = "sub-2"
$ws.Range("C4").Value = "Not Defined"
$ws.Range("E4").Value = "Experimental"
$ws.Range("F4").Value = "DCE-MRI Contrast Image sam-1"
$ws.Range("G4").Value = "Breast"
$ws.Range("H4").Value = "Not Defined"
$ws.Range("I4").Value = "Not Defined"

# Row 5: sam-2 / sub-2
$ws.Range("A5").Value = "sam-2"
$ws.Range("B5").Value = "sub-2"
$ws.Range("C5").Value = "Not Defined"
$ws.Range("E5").Value = "Experimental"
$ws.Range("F5").Value = "DCE-MRI Contrast Image sam-2"
$ws.Range("G5").Value = "Breast"
$ws.Range("H5").Value = "Not Defined"
$ws.Range("I5").Value = "Not Defined"
